# referral: changing the lab test result to CD4 lab test result
#
# Touches three sheets:
#   survey   - the question "type" text loses the " or_other" suffix
#   choices  - two whole choice rows are removed (test/viral, count/unknown),
#              three remaining choice labels are re-worded, and column C is
#              widened
#   settings - the form_title value changes and column A gets an explicit
#              width

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$choices  = $wb.Worksheets.Item("choices")
$settings = $wb.Worksheets.Item("settings")

# ---------------------------------------------------------------------------
# survey sheet: row 22 holds the "test" question; its `type` column drops the
# " or_other" suffix (becomes "select_one test " with a trailing space).
# ---------------------------------------------------------------------------
$survey.Range("A22").Value = "select_one test "

# ---------------------------------------------------------------------------
# choices sheet
# ---------------------------------------------------------------------------

# Delete the row 10 (list_name=count, name=unknown) first so row indices
# above it are untouched while we still need them.
$choices.Range("A10").EntireRow.Delete()

# Delete the row 3 (list_name=test, name=viral) choice entirely.
$choices.Range("A3").EntireRow.Delete()

# Re-word the remaining choice labels (rows have shifted up by one row for
# everything that was below row 3, but the "un"/"inconclusive"/"snooze1"/
# "snooze2" rows keep their relative order).
# After deleting row 3 ("test"/"viral"), the row that was row 6
# (result/un/Unknown) is now row 5.
$choices.Range("C5").Value = "Unknown (Lab test did not give a result after 14 days. Close this Task and set up new lab test appointment for patient to get a new lab draw)"

# The row that was row 9 (count/inconclusive/...) is now row 8.
$choices.Range("C8").Value = "Inconclusive "

# The row that was row 13 (snooze/snooze1/...) is now row 11 (after both
# deletions).
$choices.Range("C11").Value = "Okay, I will schedule a CD4 Lab Appointment`n"

# The row that was row 14 (snooze/snooze2/...) is now row 12.
$choices.Range("C12").Value = "Keep this reminder in my Task List "

# Widen column C on the choices sheet.
$choices.Columns.Item(3).ColumnWidth = 103.04666666666667

# ---------------------------------------------------------------------------
# settings sheet
# ---------------------------------------------------------------------------

# form_title value (row 2, column A).
$settings.Range("A2").Value = "CD4 Lab Test Result"

# New explicit width for column A.
$settings.Columns.Item(1).ColumnWidth = 23.046666666666667
